$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the 2000-2009 rows (old rows 2-11); everything below shifts up by 10.
$ws.Rows("2:11").Delete()

# After the shift, year 2020 (previously row 22) now lives in row 12.
# Refresh it with the higher-precision figures from the new source data.
$ws.Range("C12").Value = 10.3948359320459
$ws.Range("D12").Value = 46.2822463104876
$ws.Range("E12").Value = 43.3229177574657

# Append the two new years, copying the label-cell formatting from row 12.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A14").PasteSpecial(-4122)

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 100
$ws.Range("C13").Value = 6.4
$ws.Range("D13").Value = 54.7
$ws.Range("E13").Value = 38.9

$ws.Range("A14").Value = "2022年"
$ws.Range("B14").Value = 100
$ws.Range("C14").Value = 10.5
$ws.Range("D14").Value = 41.8
$ws.Range("E14").Value = 47.7
